$p = $ppt.ActivePresentation

# Slide 9: "What steps are we taking as Australians? – Script - XXXX" -> Natalie
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "What steps are we taking as Australians? – Script - Natalie"

# Slide 12: "So what is being done? – Script - XXXX" -> Shane
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "So what is being done? – Script - Shane"

# Slide 15: "Why is the team doing this? – Script - XXXX" -> Emily
$s15 = $p.Slides.Item(15)
$s15.Shapes.Item(1).TextFrame.TextRange.Text = "Why is the team doing this? – Script - Emily"

# Slide 18: "How can I find out more? – Script - XXXX" -> Emily
# and append " etc" to the launch sentence
$s18 = $p.Slides.Item(18)
$s18.Shapes.Item(1).TextFrame.TextRange.Text = "How can I find out more? – Script - Emily"
$s18.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "We’ll be launching soon, keep an eye out for our website launch etc"
